# Apply cryptocurrency price/volume updates to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new value looks like a plain
# number (e.g. "0.4517", "327.00"). These must stay text - exactly as
# authored - and not be re-interpreted/rounded as floating point numbers.
$textCells = @('D5', 'D7', 'D8', 'D9', 'D11', 'D12', 'D14', 'D15', 'D18', 'D19', 'D20', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '28.254.06'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '1.794.75'
$ws.Range('E3').Value = '  +1.87%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '327.00'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = '0.4517'
$ws.Range('E7').Value = '  +14.91%  '
$ws.Range('D8').Value = '0.3745'
$ws.Range('E8').Value = '  +10.36%  '
$ws.Range('D9').Value = '44.69'
$ws.Range('E9').Value = '  -1.39%  '
$ws.Range('E10').Value = '  +2.23%  '
$ws.Range('D11').Value = '0.07539'
$ws.Range('E11').Value = '  +4.27%  '
$ws.Range('D12').Value = '22.56'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').Value = '6.288'
$ws.Range('E14').Value = '  +2.27%  '
$ws.Range('D15').Value = '7.539'
$ws.Range('E15').Value = '  +6.09%  '
$ws.Range('D16').Value = '1.787.71'
$ws.Range('E16').Value = '  +5.47%  '
$ws.Range('E17').Value = '  +2.91%  '
$ws.Range('D18').Value = '0.06724'
$ws.Range('E18').Value = '  +1.47%  '
$ws.Range('D19').Value = '81.08'
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').Value = '0.9999'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('E21').Value = '  +3.56%  '
$ws.Range('D22').Value = '6.341'
$ws.Range('E22').Value = '  +1.84%  '
$ws.Range('D23').Value = '28.245.32'
$ws.Range('E23').Value = '  +0.87%  '
$ws.Range('D24').Value = '11.76'
$ws.Range('E24').Value = '  +0.83%  '
$ws.Range('D25').Value = '2.420'
$ws.Range('E25').Value = '  +1.44%  '
$ws.Range('D26').Value = '20.50'
$ws.Range('E26').Value = '  +2.77%  '
$ws.Range('D27').Value = '151.82'
$ws.Range('E27').Value = '  -1.98%  '
$ws.Range('D28').Value = '2.353'
$ws.Range('E28').Value = '  +1.68%  '
$ws.Range('D29').Value = '1.986.91'
$ws.Range('E29').Value = '  +1.35%  '
$ws.Range('D30').Value = '132.91'
$ws.Range('E30').Value = '  +2.76%  '
$ws.Range('D31').Value = '1.230'
$ws.Range('E31').Value = '  -3.80%  '
$ws.Range('D32').Value = '4.024'
$ws.Range('E32').Value = '  -1.38%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '5.817'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = '0.09421'
$ws.Range('E34').Value = '  +8.01%  '
$ws.Range('D35').Value = '0.2323'
$ws.Range('E35').Value = '  +9.95%  '
$ws.Range('D36').Value = '12.11'
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('D37').Value = '0.06337'
$ws.Range('E37').Value = '  +2.51%  '
$ws.Range('D38').Value = '0.02327'
$ws.Range('E38').Value = '  +1.57%  '
$ws.Range('D39').Value = '5.165'
$ws.Range('E39').Value = '  +0.43%  '
$ws.Range('D40').Value = '0.6560'
$ws.Range('E40').Value = '  +0.89%  '
$ws.Range('D41').Value = '8.304'
$ws.Range('E41').Value = '  +5.28%  '
$ws.Range('D42').Value = '1.470'
$ws.Range('E42').Value = '  -1.86%  '
$ws.Range('D43').Value = '1.207'
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '14.12'
$ws.Range('E44').Value = '  +2.64%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').Value = '0.9995'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '0.6094'
$ws.Range('E46').Value = '  +1.68%  '
$ws.Range('D47').Value = '3.790'
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('D48').Value = '129.91'
$ws.Range('E48').Value = '  +2.51%  '
$ws.Range('D49').Value = '2.025'
$ws.Range('E49').Value = '  +1.82%  '
$ws.Range('D50').Value = '0.07126'
$ws.Range('E50').Value = '  +1.85%  '
$ws.Range('D51').Value = '1.161'
$ws.Range('E51').Value = '  +0.27%  '
